# Tutorial 6 solution update:
#  - Reformat the Date column (A3:A21) from dd/mm/yyyy to dd-mm-yyyy
#  - Update the Total(D)/Real(E)/Absent(H)/Invalid(G) attendance counters
#    for a handful of rows to reflect the corrected attendance data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New (dashed) date strings for rows 3..21, column A
$dates = @{
    3  = "28-07-2022"
    4  = "01-08-2022"
    5  = "04-08-2022"
    6  = "08-08-2022"
    7  = "11-08-2022"
    8  = "15-08-2022"
    9  = "18-08-2022"
    10 = "22-08-2022"
    11 = "25-08-2022"
    12 = "29-08-2022"
    13 = "01-09-2022"
    14 = "05-09-2022"
    15 = "08-09-2022"
    16 = "12-09-2022"
    17 = "15-09-2022"
    18 = "19-09-2022"
    19 = "22-09-2022"
    20 = "26-09-2022"
    21 = "29-09-2022"
}

# Force the cells to keep a plain text value - otherwise Excel's automatic
# date recognition would silently turn e.g. "01-08-2022" into a real date
# serial number instead of the literal dashed string used in the source diff.
foreach ($r in $dates.Keys) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.NumberFormat = "@"
    $cell.Value = $dates[$r]
}

# Updated attendance counters: D=Total, E=Real, F=Duplicate, G=Invalid, H=Absent
$ws.Cells.Item(3, 4).Value = 1   # D3
$ws.Cells.Item(3, 7).Value = 1   # G3

$ws.Cells.Item(4, 4).Value = 1   # D4
$ws.Cells.Item(4, 5).Value = 1   # E4
$ws.Cells.Item(4, 8).Value = 0   # H4

$ws.Cells.Item(6, 4).Value = 1   # D6
$ws.Cells.Item(6, 5).Value = 1   # E6
$ws.Cells.Item(6, 8).Value = 0   # H6

$ws.Cells.Item(12, 4).Value = 1  # D12
$ws.Cells.Item(12, 5).Value = 1  # E12
$ws.Cells.Item(12, 8).Value = 0  # H12

$ws.Cells.Item(13, 4).Value = 1  # D13
$ws.Cells.Item(13, 5).Value = 1  # E13
$ws.Cells.Item(13, 8).Value = 0  # H13
